# Auto-generated edit script for EquationRelatedReactant.xlsx
# Adds product cells (B..F) for many reaction rows on sheet "生成物" (sheet 4)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

$ws.Cells.Item(314, 2).Value = "二氧化碳"
$ws.Cells.Item(314, 3).Value = "二氧化氮"
$ws.Cells.Item(314, 4).Value = "水"
$ws.Cells.Item(315, 2).Value = "二氧化碳"
$ws.Cells.Item(315, 3).Value = "二氧化硫"
$ws.Cells.Item(315, 4).Value = "水"
$ws.Cells.Item(316, 2).Value = "碳酸氢根离子"
$ws.Cells.Item(316, 3).Value = "氢氧根离子"
$ws.Cells.Item(317, 2).Value = "硅酸钙"
$ws.Cells.Item(317, 3).Value = "二氧化碳"
$ws.Cells.Item(318, 2).Value = "硅酸钠"
$ws.Cells.Item(318, 3).Value = "二氧化碳"
$ws.Cells.Item(321, 2).Value = "四氯化硅"
$ws.Cells.Item(321, 3).Value = "硅单质"
$ws.Cells.Item(321, 4).Value = "氯化氢"
$ws.Cells.Item(323, 2).Value = "碳酸钠"
$ws.Cells.Item(323, 3).Value = "氧气"
$ws.Cells.Item(324, 2).Value = "氢氧化钠"
$ws.Cells.Item(324, 3).Value = "氧气"
$ws.Cells.Item(325, 2).Value = "氯化钠"
$ws.Cells.Item(325, 3).Value = "水"
$ws.Cells.Item(325, 4).Value = "氧气"
$ws.Cells.Item(326, 2).Value = "水"
$ws.Cells.Item(327, 2).Value = "氢氧化铜"
$ws.Cells.Item(327, 3).Value = "硫酸钠"
$ws.Cells.Item(327, 4).Value = "氢气"
$ws.Cells.Item(328, 2).Value = "硫酸铁"
$ws.Cells.Item(328, 3).Value = "二氧化硫"
$ws.Cells.Item(328, 4).Value = "水"
$ws.Cells.Item(329, 2).Value = "硝酸铁"
$ws.Cells.Item(329, 3).Value = "水"
$ws.Cells.Item(329, 4).Value = "一氧化氮"
$ws.Cells.Item(330, 2).Value = "亚铁离子"
$ws.Cells.Item(331, 2).Value = "亚铁离子"
$ws.Cells.Item(331, 3).Value = "银单质"
$ws.Cells.Item(332, 2).Value = "亚铁离子"
$ws.Cells.Item(332, 3).Value = "锌单质"
$ws.Cells.Item(333, 2).Value = "氢氧化铁"
$ws.Cells.Item(333, 3).Value = "铵根离子"
$ws.Cells.Item(334, 2).Value = "铅单质"
$ws.Cells.Item(334, 3).Value = "二氧化铅"
$ws.Cells.Item(334, 4).Value = "硫酸"
$ws.Cells.Item(335, 2).Value = "硫酸铅"
$ws.Cells.Item(335, 3).Value = "水"
$ws.Cells.Item(336, 2).Value = "铜离子"
$ws.Cells.Item(336, 3).Value = "汞单质"
$ws.Cells.Item(337, 2).Value = "硝酸铜"
$ws.Cells.Item(337, 3).Value = "二氧化氮"
$ws.Cells.Item(337, 4).Value = "水"
$ws.Cells.Item(338, 2).Value = "硫酸铜"
$ws.Cells.Item(338, 3).Value = "二氧化硫"
$ws.Cells.Item(338, 4).Value = "水"
$ws.Cells.Item(339, 2).Value = "硝酸铜"
$ws.Cells.Item(339, 3).Value = "一氧化氮"
$ws.Cells.Item(339, 4).Value = "水"
$ws.Cells.Item(340, 2).Value = "亚铁离子"
$ws.Cells.Item(340, 3).Value = "铜离子"
$ws.Cells.Item(341, 2).Value = "硫化铜"
$ws.Cells.Item(342, 2).Value = "偏铝酸根"
$ws.Cells.Item(342, 3).Value = "氢气"
$ws.Cells.Item(343, 2).Value = "氢氧化铝"
$ws.Cells.Item(343, 3).Value = "二氧化碳"
$ws.Cells.Item(343, 4).Value = "水"
$ws.Cells.Item(344, 2).Value = "氨气"
$ws.Cells.Item(344, 3).Value = "水"
$ws.Cells.Item(345, 2).Value = "氨气"
$ws.Cells.Item(345, 3).Value = "氢气"
$ws.Cells.Item(345, 4).Value = "镁离子"
$ws.Cells.Item(346, 2).Value = "银氨离子"
$ws.Cells.Item(347, 2).Value = "氯化银"
$ws.Cells.Item(348, 2).Value = "氯化铝"
$ws.Cells.Item(348, 3).Value = "亚硫酸铝"
$ws.Cells.Item(348, 4).Value = "硫单质"
$ws.Cells.Item(349, 2).Value = "银单质"
$ws.Cells.Item(349, 3).Value = "氢氧化锌"
$ws.Cells.Item(350, 2).Value = "锰离子"
$ws.Cells.Item(350, 3).Value = "氧气"
$ws.Cells.Item(350, 4).Value = "水"
$ws.Cells.Item(351, 2).Value = "2-丙醇"
$ws.Cells.Item(352, 2).Value = "乙烷"
$ws.Cells.Item(353, 2).Value = "氯乙烯"
$ws.Cells.Item(354, 2).Value = "乙醛"
$ws.Cells.Item(355, 2).Value = "1,1,2,2-四溴乙烷"
$ws.Cells.Item(356, 2).Value = "乙烷"
$ws.Cells.Item(357, 2).Value = "一氯乙烷"
$ws.Cells.Item(358, 2).Value = "1，2-二氯乙烷"
$ws.Cells.Item(359, 2).Value = "乙醇"
$ws.Cells.Item(360, 2).Value = "1,2-二溴乙烷"
$ws.Cells.Item(361, 2).Value = "乙醇"
$ws.Cells.Item(363, 2).Value = "环己烷"
$ws.Cells.Item(364, 2).Value = "环己醇"
$ws.Cells.Item(367, 2).Value = "聚丙烯"
$ws.Cells.Item(368, 2).Value = "聚丙烯腈"
$ws.Cells.Item(369, 2).Value = "聚乙烯"
$ws.Cells.Item(370, 2).Value = "聚四氟乙烯"
$ws.Cells.Item(371, 2).Value = "聚异戊二烯"
$ws.Cells.Item(372, 2).Value = "聚苯乙烯"
$cell = $ws.Cells.Item(373, 2)
$cell.Value = "聚1，3-丁二烯"
$rc2 = $cell.Characters(3,1)
$rc2.Font.Name = "Droid Sans Fallback"
$rc2.Font.Size = 10
$rc2.Font.Color = 0
$rc3 = $cell.Characters(4,2)
$rc3.Font.Name = "Arial"
$rc3.Font.Size = 10
$rc3.Font.Color = 0
$rc4 = $cell.Characters(6,3)
$rc4.Font.Name = "Droid Sans Fallback"
$rc4.Font.Size = 10
$rc4.Font.Color = 0
$ws.Cells.Item(375, 2).Value = "乙酸乙酯"
$ws.Cells.Item(375, 3).Value = "水"
$ws.Cells.Item(377, 2).Value = "乙酸乙酯"
$ws.Cells.Item(377, 3).Value = "水"
$ws.Cells.Item(378, 2).Value = "溴乙烷"
$ws.Cells.Item(378, 3).Value = "水"
$ws.Cells.Item(379, 2).Value = "乙烯"
$ws.Cells.Item(379, 3).Value = "水"
$ws.Cells.Item(380, 2).Value = "乙醚"
$ws.Cells.Item(380, 3).Value = "水"
$ws.Cells.Item(381, 2).Value = "乙醇"
$ws.Cells.Item(381, 3).Value = "溴化钠"
$ws.Cells.Item(382, 2).Value = "一氯甲烷"
$ws.Cells.Item(382, 3).Value = "二氯甲烷"
$ws.Cells.Item(382, 4).Value = "三氯甲烷"
$ws.Cells.Item(382, 5).Value = "四氯化碳"
$ws.Cells.Item(382, 6).Value = "氯化氢"
$ws.Cells.Item(383, 2).Value = "2-硝基甲苯"
$ws.Cells.Item(383, 3).Value = "4-硝基甲苯"
$ws.Cells.Item(383, 4).Value = "水"
$ws.Cells.Item(384, 2).Value = "2,4,6-三硝基甲苯"
$ws.Cells.Item(384, 3).Value = "水"
$ws.Cells.Item(386, 2).Value = "溴苯"
$ws.Cells.Item(386, 3).Value = "溴化氢"
$ws.Cells.Item(387, 2).Value = "溴苯"
$ws.Cells.Item(387, 3).Value = "溴化氢"
$ws.Cells.Item(388, 2).Value = "硝基苯"
$ws.Cells.Item(388, 3).Value = "水"
$ws.Cells.Item(389, 2).Value = "三硝基苯酚"
$ws.Cells.Item(389, 3).Value = "水"
$ws.Cells.Item(390, 2).Value = "2,4,6-三溴苯酚"
$ws.Cells.Item(390, 3).Value = "溴化氢"
$ws.Cells.Item(391, 2).Value = "乙烯"
$ws.Cells.Item(391, 3).Value = "水"
$ws.Cells.Item(392, 2).Value = "乙醇"
$ws.Cells.Item(392, 3).Value = "溴化钠"
$ws.Cells.Item(392, 4).Value = "水"
$ws.Cells.Item(393, 2).Value = "酚醛树脂"
$ws.Cells.Item(393, 3).Value = "水"
$ws.Cells.Item(394, 2).Value = "二氧化碳"
$ws.Cells.Item(394, 3).Value = "水"
$ws.Cells.Item(395, 2).Value = "乙醛"

# Add a new (incremented) hidden-filter-database defined name for the "反应物" sheet,
# mirroring the pattern of previous _FilterDatabase_0... entries already present.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Names.Add("_xlnm._FilterDatabase_0_0_0_0", "=反应物!`$A`$1:`$C`$853")

# Update the view/selection on the "生成物" sheet to reflect where the user ended up editing.
$ws.Activate()
$ws.Range("B396").Select()
$wb.Application.ActiveWindow.ScrollRow = 388
$wb.Application.ActiveWindow.ScrollColumn = 1

